$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.478793263435364
$ws.Range("B1").Value = 1.898364901542664
$ws.Range("C1").Value = 2.849325180053711
$ws.Range("D1").Value = 5.278578758239746
$ws.Range("E1").Value = 1.123391270637512
